# Apply updates to Bulgaria_FX.xlsx per commit diff:
#  - Row 257: correct D257 and F257 values
#  - Append rows 258-260 with new monthly FX data points

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix existing row 257 values ---
$ws.Range("D257").Value = 1.81045
$ws.Range("F257").Value = 1.79962

# --- Row 258 ---
$ws.Range("A258").Value = 45170.33333333334
$ws.Range("A258").Style = $ws.Range("A257").Style
$ws.Range("B258").Value = "FX_IDC:USDBGN"
$ws.Range("C258").Value = 1.80361
$ws.Range("D258").Value = 1.85633
$ws.Range("E258").Value = 1.80361
$ws.Range("F258").Value = 1.84617
$ws.Range("G258").Value = 0

# --- Row 259 ---
$ws.Range("A259").Value = 45201.375
$ws.Range("A259").Style = $ws.Range("A257").Style
$ws.Range("B259").Value = "FX_IDC:USDBGN"
$ws.Range("C259").Value = 1.85739
$ws.Range("D259").Value = 1.86821
$ws.Range("E259").Value = 1.83957
$ws.Range("F259").Value = 1.84182
$ws.Range("G259").Value = 0

# --- Row 260 ---
$ws.Range("A260").Value = 45231.375
$ws.Range("A260").Style = $ws.Range("A257").Style
$ws.Range("B260").Value = "FX_IDC:USDBGN"
$ws.Range("C260").Value = 1.85615
$ws.Range("D260").Value = 1.85615
$ws.Range("E260").Value = 1.8209
$ws.Range("F260").Value = 1.83302
$ws.Range("G260").Value = 0
